# Update the "two-digit-mul" worksheet: bump the date and replace each
# multiplication problem's operands with the new values from the commit.

$d = $word.ActiveDocument

# Map of old text -> new text, in document order.
$replacements = [ordered]@{
    "2023-11-27 Monday" = "2023-11-28 Tuesday"
    "94×88="            = "68×62="
    "33×95="            = "12×94="
    "38×47="            = "43×23="
    "16×21="            = "19×40="
    "69×37="            = "60×16="
    "66×43="            = "74×37="
    "35×87="            = "34×79="
    "55×28="            = "48×85="
    "52×57="            = "75×41="
    "86×30="            = "31×97="
    "65×21="            = "43×75="
    "11×37="            = "90×73="
    "49×83="            = "77×79="
    "49×91="            = "68×55="
    "50×89="            = "81×38="
    "29×71="            = "89×91="
    "96×98="            = "63×15="
    "94×12="            = "87×87="
    "40×41="            = "89×61="
    "68×38="            = "94×70="
    "26×59="            = "79×23="
    "72×65="            = "64×58="
    "30×30="            = "74×45="
    "64×12="            = "67×88="
    "35×43="            = "14×66="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
